$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: capture the existing 30 data rows (rows 2-31, columns A-E) into memory.
# Column A holds dates; read them as their underlying numeric (OADate) serial so
# that writing them back later does not trigger auto date-formatting on cells
# that don't already carry the date style.
$existing = @()
for ($r = 2; $r -le 31; $r++) {
    $row = @()
    for ($c = 1; $c -le 5; $c++) {
        $v = $ws.Cells.Item($r, $c).Value()
        if ($c -eq 1) {
            $v = $v.ToOADate()
        }
        $row += $v
    }
    $existing += , $row
}

# Step 2: write the captured rows back, shifted down by 11 rows (new rows 13-42)
for ($i = 0; $i -lt $existing.Count; $i++) {
    $r = 13 + $i
    $row = $existing[$i]
    for ($c = 1; $c -le 5; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}

# Step 3: new data for the 11 inserted rows (years 1985-1995), rows 2-12
$newData = @(
    @(31228, 1985, 1.777145513200118,   1986, 1.897195788720762),
    @(31593, 1986, 1.520629350269975,   1987, 1.941216181406236),
    @(31958, 1987, 0.006150708382479664,1988, 2.273434064658209),
    @(32324, 1988, 2.284628917872622,   1989, 2.270390965314983),
    @(32689, 1989, 3.661580277249166,   1990, 2.193258610001214),
    @(33054, 1990, 4.130619852766437,   1991, 2.272904440822465),
    @(33419, 1991, 6.364491101711689,   1992, 2.793798186209284),
    @(33785, 1992, 2.932796654414149,   1993, 2.584450468619459),
    @(34150, 1993, -1.06363680093724,   1994, 2.368493192930488),
    @(34515, 1994, 2.479893153134016,   1995, 2.567096653116252),
    @(34880, 1995, 2.432437183852798,   1996, 2.872765583543457)
)

for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = 2 + $i
    $row = $newData[$i]
    for ($c = 1; $c -le 5; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}

# Step 4: copy the date number format/style from an already-correctly-styled
# date cell onto every date cell in column A (A2:A42), so they all carry the
# same style index instead of a freshly minted / default one. This covers the
# newly inserted rows (2-12) as well as the rows that moved into previously
# nonexistent rows (32-42).
$ws.Range("A13").Copy()
$ws.Range("A2:A42").PasteSpecial(-4122)
$excel.CutCopyMode = 0

Write-Output "edit complete"
Write-Output $ws.UsedRange.Address()
